$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45180 -> 2023-09-11)
# that was bumped by one day (45181 -> 2023-09-12) for every data row
# (rows 2 through 528).
for ($r = 2; $r -le 528; $r++) {
    $ws.Cells.Item($r, 3).Value = 45181
}
